# price fix + incredible
# Update four Digikala product links to the shorter "fresh" product URL format.
# Order matters: Excel appends newly-introduced shared strings in the order
# they are first written, and the target workbook appends them in this order:
#   dkp-4714424, dkp-757476, dkp-1023378, dkp-1485597
# which corresponds to rows 40, 31, 9, 10 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "https://www.digikala.com/fresh/product/dkp-4714424/"
$ws.Range("A31").Value = "https://www.digikala.com/fresh/product/dkp-757476/"
$ws.Range("A9").Value  = "https://www.digikala.com/fresh/product/dkp-1023378/"
$ws.Range("A10").Value = "https://www.digikala.com/fresh/product/dkp-1485597/"

# Move the active selection from A12 to A9.
$ws.Range("A9").Select()
